$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("CCDeferredCC_27").Range("B2").Value = "Sat Nov 08 13:43:08 IST 2025"
$wb.Worksheets.Item("CMCAutopayCC_27").Range("B2").Value = "Sat Nov 08 13:49:13 IST 2025"
$wb.Worksheets.Item("PayNowDCFCC_27").Range("B2").Value = "Sat Nov 08 13:53:29 IST 2025"
